# Insert a new column at AI (35th column). This shifts the previous
# AI column (NPCType) to AJ, and the previous AJ column (DescID) to AK,
# copying formatting along the way (matches native Excel "Insert" behavior).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("AI").Insert()

# The newly inserted column takes on the same display width as its
# neighbour (column AH), mirroring the target layout.
$ws.Columns("AI").ColumnWidth = $ws.Columns("AH").ColumnWidth

# Populate the header / metadata rows (1-10) for the newly inserted column.
$ws.Range("AI1").Value2  = "AIOwnerID"
$ws.Range("AI2").Value2  = "object"
$ws.Range("AI3").Value2  = 1
$ws.Range("AI4").Value2  = 1
$ws.Range("AI5").Value2  = 1
$ws.Range("AI6").Value2  = 0
$ws.Range("AI7").Value2  = 0
$ws.Range("AI8").Value2  = 0
$ws.Range("AI9").Value2  = 0
$ws.Range("AI10").Value2 = "AI"

# Populate the data rows (11-62) for the new column with the default value.
for ($r = 11; $r -le 62; $r++) {
    $ws.Range("AI$r").Value2 = 0
}

# The NPCType column (shifted from AI to AJ) is also flipped to
# Public=TRUE / Private=TRUE as part of this change.
$ws.Range("AJ3").Value2 = 1
$ws.Range("AJ4").Value2 = 1

# Restore the active selection to match the edited workbook.
$ws.Range("AI11").Select()
